$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Switch" row (row 6): new footprint link for the slide switch ---
$ws.Range("B6").Value = "https://www.lcsc.com/product-detail/span-style-background-color-ff0-Slide-span-Switches_SHOU-HAN-MSK12CO2-SZ_C2681568.html"

# --- Update the "CR2032 Battery holder" row (row 3): new footprint link for the battery connector ---
$ws.Range("B3").Value = "https://www.lcsc.com/product-detail/span-style-background-color-ff0-Battery-span-Connectors_Q-J-C70373_C70373.html"

# Row 3 (battery holder link) did not previously carry a hyperlink/Link style - add it now,
# matching the existing hyperlinked style used for the Switch row (B6).
$ws.Hyperlinks.Add($ws.Range("B3"), $ws.Range("B3").Value())
$ws.Range("B3").Style = "Link"

# --- Update the manually-set selection to reflect the last edited cell ---
$ws.Range("B3").Select()

Write-Host "Materialliste updated"
